# Add "TransactionType16" column to the RawEarnings sheet (insert before the
# existing last column, "ApprenticeshipContractType"), populate it with 0s,
# and make RawEarnings the active/selected sheet with AB1 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RawEarnings")

# Insert a new column at AB (pushes the old AB "ApprenticeshipContractType"
# column, and its data, to AC).
$ws.Columns("AB:AB").Insert()

# Match the width used by the neighbouring TransactionTypeNN columns.
$ws.Columns("AB:AB").ColumnWidth = 16.83

# New header + zero-filled data for the new TransactionType16 column.
$ws.Range("AB1").Value = "TransactionType16"
for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 28).Value = 0
}

# Make RawEarnings the active sheet with AB1 selected (matches the saved view
# state in the workbook).
$ws.Activate()
$ws.Range("AB1").Select()
